$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Expand used range: dimension goes from A1:H21 to A1:H31 (10 new rows of data).
# The timestamp (A) / label (B) columns are untouched for existing rows 2-21 and simply
# extended linearly (+100 per row) for the 10 new rows (22-31).
# The sensor columns (C:H) are rewritten: 3 new rows of data are inserted logically at the
# top of the numeric series (now occupying rows 2-4), the old values for rows 2-18 shift down
# to rows 5-21, the old values for rows 19-21 shift down to rows 22-24, and 7 more new rows of
# data are appended at rows 25-31.

# New sensor data (C:H) for the first 3 new rows (2-4), A/B already correct (0/100/200, walkingToRunning)
$ws.Cells.Item(2, 3).Value = 1.45181941986084
$ws.Cells.Item(2, 4).Value = -10.64814472198486
$ws.Cells.Item(2, 5).Value = 5.583051204681396
$ws.Cells.Item(2, 6).Value = -0.7805059552192688
$ws.Cells.Item(2, 7).Value = -0.7913583517074585
$ws.Cells.Item(2, 8).Value = 0.3286340832710266
$ws.Cells.Item(3, 3).Value = 12.47824478149414
$ws.Cells.Item(3, 4).Value = -36.90848541259766
$ws.Cells.Item(3, 5).Value = 23.18166542053223
$ws.Cells.Item(3, 6).Value = 0.1883520781993866
$ws.Cells.Item(3, 7).Value = 0.3149188160896301
$ws.Cells.Item(3, 8).Value = -0.798948347568512
$ws.Cells.Item(4, 3).Value = 4.87645149230957
$ws.Cells.Item(4, 4).Value = -11.23874282836914
$ws.Cells.Item(4, 5).Value = 11.8306131362915
$ws.Cells.Item(4, 6).Value = 2.627274990081787
$ws.Cells.Item(4, 7).Value = -2.136787414550781
$ws.Cells.Item(4, 8).Value = -1.515871286392212

# Shift old sensor data (rows 2-18) down into rows 5-21 (A/B columns already correct there)
$ws.Cells.Item(5, 3).Value = 3.935386657714844
$ws.Cells.Item(5, 4).Value = -20.30976104736328
$ws.Cells.Item(5, 5).Value = 10.9278678894043
$ws.Cells.Item(5, 6).Value = -0.4161854982376098
$ws.Cells.Item(5, 7).Value = 2.661696434020996
$ws.Cells.Item(5, 8).Value = -1.426389098167419
$ws.Cells.Item(6, 3).Value = 6.922998905181885
$ws.Cells.Item(6, 4).Value = -22.64658546447754
$ws.Cells.Item(6, 5).Value = 13.07839107513428
$ws.Cells.Item(6, 6).Value = -3.087868690490723
$ws.Cells.Item(6, 7).Value = -3.132277011871338
$ws.Cells.Item(6, 8).Value = 3.167563915252685
$ws.Cells.Item(7, 3).Value = -17.48348999023438
$ws.Cells.Item(7, 4).Value = 3.01125955581665
$ws.Cells.Item(7, 5).Value = -5.847013473510742
$ws.Cells.Item(7, 6).Value = -8.064251899719238
$ws.Cells.Item(7, 7).Value = 4.051333904266357
$ws.Cells.Item(7, 8).Value = 0.5699164867401123
$ws.Cells.Item(8, 3).Value = -9.65645980834961
$ws.Cells.Item(8, 4).Value = -6.639880180358887
$ws.Cells.Item(8, 5).Value = 1.525803089141846
$ws.Cells.Item(8, 6).Value = 3.350589513778687
$ws.Cells.Item(8, 7).Value = -4.077699184417725
$ws.Cells.Item(8, 8).Value = -0.3270362019538879
$ws.Cells.Item(9, 3).Value = 20.44024848937988
$ws.Cells.Item(9, 4).Value = -40.82086944580078
$ws.Cells.Item(9, 5).Value = 14.10586166381836
$ws.Cells.Item(9, 6).Value = 2.837132215499878
$ws.Cells.Item(9, 7).Value = -3.239868640899658
$ws.Cells.Item(9, 8).Value = -1.359277486801148
$ws.Cells.Item(10, 3).Value = -15.60717296600342
$ws.Cells.Item(10, 4).Value = -21.49799346923828
$ws.Cells.Item(10, 5).Value = 26.6249885559082
$ws.Cells.Item(10, 6).Value = 6.7850022315979
$ws.Cells.Item(10, 7).Value = 5.038301467895508
$ws.Cells.Item(10, 8).Value = 0.312122493982315
$ws.Cells.Item(11, 3).Value = 1.734474897384644
$ws.Cells.Item(11, 4).Value = 2.406877994537353
$ws.Cells.Item(11, 5).Value = 10.91264820098877
$ws.Cells.Item(11, 6).Value = 2.806239366531372
$ws.Cells.Item(11, 7).Value = 1.31573474407196
$ws.Cells.Item(11, 8).Value = -1.794970631599426
$ws.Cells.Item(12, 3).Value = 31.86017036437988
$ws.Cells.Item(12, 4).Value = -35.68617630004883
$ws.Cells.Item(12, 5).Value = 40.4693717956543
$ws.Cells.Item(12, 6).Value = -4.413057804107666
$ws.Cells.Item(12, 7).Value = 0.4736432135105133
$ws.Cells.Item(12, 8).Value = 2.279133319854736
$ws.Cells.Item(13, 3).Value = -57.41475677490234
$ws.Cells.Item(13, 4).Value = 17.37560653686523
$ws.Cells.Item(13, 5).Value = -13.54984855651856
$ws.Cells.Item(13, 6).Value = -5.49536657333374
$ws.Cells.Item(13, 7).Value = -2.257694959640503
$ws.Cells.Item(13, 8).Value = -2.082059383392334
$ws.Cells.Item(14, 3).Value = -13.51763439178467
$ws.Cells.Item(14, 4).Value = -10.28693199157715
$ws.Cells.Item(14, 5).Value = -0.659794807434082
$ws.Cells.Item(14, 6).Value = 0.6288388967514038
$ws.Cells.Item(14, 7).Value = 10.04877281188965
$ws.Cells.Item(14, 8).Value = -1.388039588928223
$ws.Cells.Item(15, 3).Value = -20.15132904052734
$ws.Cells.Item(15, 4).Value = -55.33749389648438
$ws.Cells.Item(15, 5).Value = 21.29398727416992
$ws.Cells.Item(15, 6).Value = 4.64808177947998
$ws.Cells.Item(15, 7).Value = -3.816709280014038
$ws.Cells.Item(15, 8).Value = 2.034122467041016
$ws.Cells.Item(16, 3).Value = 15.73336029052734
$ws.Cells.Item(16, 4).Value = -4.189780712127686
$ws.Cells.Item(16, 5).Value = 11.97433757781982
$ws.Cells.Item(16, 6).Value = 6.68113899230957
$ws.Cells.Item(16, 7).Value = 7.193863868713379
$ws.Cells.Item(16, 8).Value = 0.3451456725597381
$ws.Cells.Item(17, 3).Value = 15.30455780029297
$ws.Cells.Item(17, 4).Value = -12.82610034942627
$ws.Cells.Item(17, 5).Value = 17.12981605529785
$ws.Cells.Item(17, 6).Value = -2.778942108154297
$ws.Cells.Item(17, 7).Value = -2.685398578643799
$ws.Cells.Item(17, 8).Value = 1.103613972663879
$ws.Cells.Item(18, 3).Value = -34.54495620727539
$ws.Cells.Item(18, 4).Value = -21.88593673706055
$ws.Cells.Item(18, 5).Value = -7.297530174255371
$ws.Cells.Item(18, 6).Value = -6.552108764648438
$ws.Cells.Item(18, 7).Value = 1.416934847831726
$ws.Cells.Item(18, 8).Value = 2.430400848388672
$ws.Cells.Item(19, 3).Value = 34.65779495239258
$ws.Cells.Item(19, 4).Value = 5.878978252410889
$ws.Cells.Item(19, 5).Value = 2.555978298187256
$ws.Cells.Item(19, 6).Value = -3.362174272537231
$ws.Cells.Item(19, 7).Value = 17.18923950195312
$ws.Cells.Item(19, 8).Value = -1.967010855674744
$ws.Cells.Item(20, 3).Value = 2.989975452423096
$ws.Cells.Item(20, 4).Value = -18.37696075439453
$ws.Cells.Item(20, 5).Value = -2.741829395294189
$ws.Cells.Item(20, 6).Value = 1.812880396842956
$ws.Cells.Item(20, 7).Value = -3.756521940231323
$ws.Cells.Item(20, 8).Value = -3.184607982635498
$ws.Cells.Item(21, 3).Value = 13.91487789154053
$ws.Cells.Item(21, 4).Value = -47.14267730712891
$ws.Cells.Item(21, 5).Value = 13.05474662780762
$ws.Cells.Item(21, 6).Value = 8.206597328186035
$ws.Cells.Item(21, 7).Value = -8.100137710571289
$ws.Cells.Item(21, 8).Value = -2.368615627288818

# Shift old sensor data (rows 19-21) down into rows 22-24; also need new A/B there since those rows did not exist before
$ws.Cells.Item(22, 1).Value = 2000
$ws.Cells.Item(22, 2).Value = "walkingToRunning"
$ws.Cells.Item(22, 3).Value = 0.2568814754486084
$ws.Cells.Item(22, 4).Value = 1.619793891906739
$ws.Cells.Item(22, 5).Value = 9.530179977416992
$ws.Cells.Item(22, 6).Value = 3.973769426345825
$ws.Cells.Item(22, 7).Value = -2.186322212219238
$ws.Cells.Item(22, 8).Value = -1.423193335533142
$ws.Cells.Item(23, 1).Value = 2100
$ws.Cells.Item(23, 2).Value = "walkingToRunning"
$ws.Cells.Item(23, 3).Value = -26.04407691955566
$ws.Cells.Item(23, 4).Value = -81.4648208618164
$ws.Cells.Item(23, 5).Value = 51.38869094848633
$ws.Cells.Item(23, 6).Value = -2.724081039428711
$ws.Cells.Item(23, 7).Value = 3.562910079956055
$ws.Cells.Item(23, 8).Value = 0.9864348769187928
$ws.Cells.Item(24, 1).Value = 2200
$ws.Cells.Item(24, 2).Value = "walkingToRunning"
$ws.Cells.Item(24, 3).Value = 7.554911613464355
$ws.Cells.Item(24, 4).Value = 0.749316930770874
$ws.Cells.Item(24, 5).Value = -22.46823501586914
$ws.Cells.Item(24, 6).Value = -11.11730003356934
$ws.Cells.Item(24, 7).Value = -16.41932106018066
$ws.Cells.Item(24, 8).Value = -2.057025671005249

# Append 7 brand new rows of data at rows 25-31 (timestamps 2300-2900)
$ws.Cells.Item(25, 1).Value = 2300
$ws.Cells.Item(25, 2).Value = "walkingToRunning"
$ws.Cells.Item(25, 3).Value = -8.727315902709961
$ws.Cells.Item(25, 4).Value = -7.33059024810791
$ws.Cells.Item(25, 5).Value = 0.0287117958068847
$ws.Cells.Item(25, 6).Value = 6.472879886627197
$ws.Cells.Item(25, 7).Value = -9.945176124572754
$ws.Cells.Item(25, 8).Value = 9.769540786743164
$ws.Cells.Item(26, 1).Value = 2400
$ws.Cells.Item(26, 2).Value = "walkingToRunning"
$ws.Cells.Item(26, 3).Value = -17.04781913757324
$ws.Cells.Item(26, 4).Value = -51.86902236938477
$ws.Cells.Item(26, 5).Value = 16.61567497253418
$ws.Cells.Item(26, 6).Value = 3.245661020278931
$ws.Cells.Item(26, 7).Value = -6.337657928466797
$ws.Cells.Item(26, 8).Value = 0.1816275864839553
$ws.Cells.Item(27, 1).Value = 2500
$ws.Cells.Item(27, 2).Value = "walkingToRunning"
$ws.Cells.Item(27, 3).Value = 27.31503295898437
$ws.Cells.Item(27, 4).Value = 1.940977096557617
$ws.Cells.Item(27, 5).Value = 14.65385055541992
$ws.Cells.Item(27, 6).Value = 5.774599075317383
$ws.Cells.Item(27, 7).Value = 6.254833221435547
$ws.Cells.Item(27, 8).Value = 0.6913566589355469
$ws.Cells.Item(28, 1).Value = 2600
$ws.Cells.Item(28, 2).Value = "walkingToRunning"
$ws.Cells.Item(28, 3).Value = 3.707320690155029
$ws.Cells.Item(28, 4).Value = -6.758492469787598
$ws.Cells.Item(28, 5).Value = 5.369882583618164
$ws.Cells.Item(28, 6).Value = -2.965895891189575
$ws.Cells.Item(28, 7).Value = -1.341034770011902
$ws.Cells.Item(28, 8).Value = 1.717206358909607
$ws.Cells.Item(29, 1).Value = 2700
$ws.Cells.Item(29, 2).Value = "walkingToRunning"
$ws.Cells.Item(29, 3).Value = 18.9058609008789
$ws.Cells.Item(29, 4).Value = -66.31611633300781
$ws.Cells.Item(29, 5).Value = 27.15024948120117
$ws.Cells.Item(29, 6).Value = -4.298542022705078
$ws.Cells.Item(29, 7).Value = 1.24329674243927
$ws.Cells.Item(29, 8).Value = 0.5667206645011902
$ws.Cells.Item(30, 1).Value = 2800
$ws.Cells.Item(30, 2).Value = "walkingToRunning"
$ws.Cells.Item(30, 3).Value = -72.50922393798828
$ws.Cells.Item(30, 4).Value = 29.85196113586425
$ws.Cells.Item(30, 5).Value = -13.03144264221191
$ws.Cells.Item(30, 6).Value = -6.537195205688477
$ws.Cells.Item(30, 7).Value = -4.836699962615967
$ws.Cells.Item(30, 8).Value = -2.582733631134033
$ws.Cells.Item(31, 1).Value = 2900
$ws.Cells.Item(31, 2).Value = "walkingToRunning"
$ws.Cells.Item(31, 3).Value = -21.97713661193848
$ws.Cells.Item(31, 4).Value = -12.47497940063477
$ws.Cells.Item(31, 5).Value = 15.31039047241211
$ws.Cells.Item(31, 6).Value = -1.381847739219666
$ws.Cells.Item(31, 7).Value = 14.15004062652588
$ws.Cells.Item(31, 8).Value = 5.833388328552246
